$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# === Step 1: extend formatting (before values change) ===
# Header-row style (bold / border / centered) -> new header cells H1:J1
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)  # xlPasteFormats
# Index-column style -> new index rows A8:A14
$ws.Range("A7").Copy()
$ws.Range("A8:A14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# === Step 2: seed true/false text cells by copying existing text cells
#     (typing the literal word true/false would auto-convert to a Boolean cell) ===
# G2 currently holds the text "true"; G4 currently holds the text "false"
$ws.Range("G2").Copy()
$ws.Range("J6").PasteSpecial(-4163)  # xlPasteValues -> "true"
$ws.Range("G4").Copy()
$ws.Range("J2").PasteSpecial(-4163)  # xlPasteValues -> "false"
$ws.Range("J8").PasteSpecial(-4163)  # xlPasteValues -> "false"
$ws.Range("J12").PasteSpecial(-4163)  # xlPasteValues -> "false"
$excel.CutCopyMode = 0

# === Step 3: header row text ===
$ws.Range("B1").Value = "INDEX"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "sib.INDEX"
$ws.Range("E1").Value = "sib.name"
$ws.Range("F1").Value = "sib.ph"
$ws.Range("G1").Value = "sib.addr"
$ws.Range("H1").Value = "frnds.INDEX"
$ws.Range("I1").Value = "frnds.b"
$ws.Range("J1").Value = "frnds.best"

# === Step 4: data rows ===
# row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "abhi"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 2345
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 3
# row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 3456
$ws.Range("G3").Value = "adadaa"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = "null"
# row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = "null"
# row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 0
$ws.Range("C5").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = "null"
# row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "aditi"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 98765
$ws.Range("G6").Value = "null"
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 3
# row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 1
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 12345
$ws.Range("G7").Value = "null"
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 4
$ws.Range("J7").Value = "null"
# row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 1
$ws.Range("C8").ClearContents()
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = "null"
$ws.Range("G8").Value = "fjaslkff kjas"
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 3
# row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 1
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 4
$ws.Range("J9").Value = "null"
# row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 1
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()
$ws.Range("G10").ClearContents()
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = "null"
# row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 1
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 6
$ws.Range("J11").Value = "null"
# row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = "ashish"
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 2345
$ws.Range("G12").Value = "null"
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 9
# row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 2
$ws.Range("C13").ClearContents()
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 3456
$ws.Range("G13").Value = "kjljl"
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 10
$ws.Range("J13").Value = "null"
# row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 2
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("G14").ClearContents()
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = 11
$ws.Range("J14").Value = "null"
